# Weekly update: insert a new price observation as row 203, pushing the
# existing rows 203:232 down to 204:233 (dimension grows from R232 to R233).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(203).Insert()

$ws.Range("A203").Value = 3
$ws.Range("B203").Value = "Femacal de La Calera"
$ws.Range("C203").Value = "Coquimbo"
$ws.Range("D203").Value = 44491
$ws.Range("E203").Value = 5
$ws.Range("F203").Value = 100112032
$ws.Range("G203").Value = "Zapallo italiano"
$ws.Range("H203").Value = "Sin especificar"
$ws.Range("I203").Value = "Primera"
$ws.Range("J203").Value = 90
$ws.Range("K203").Value = 9500
$ws.Range("L203").Value = 10000
$ws.Range("M203").Value = 9750
$ws.Range("N203").Value = "$/caja 70 unidades"
$ws.Range("O203").Value = "Región de Arica y Parinacota"
$ws.Range("P203").Value = 139
$ws.Range("Q203").Value = 70
$ws.Range("R203").Value = "Hortaliza"
